# AFLFixtures2023.xlsx -- "Add files via upload"
# Refresh the Fixture sheet kickoff date/times in column C for rows 128-199
# (match rounds originally saved with placeholder round-start timestamps,
# here updated to the actual per-game kickoff date/time), and move the
# sheet selection to where the author last left off (D196, scrolled near
# row 175).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fixture")

$ws.Range("C128").Value = 45106.805555555598
$ws.Range("C129").Value = 45107.826388888898
$ws.Range("C130").Value = 45108.572916666701
$ws.Range("C131").Value = 45108.572916666701
$ws.Range("C132").Value = 45108.690972222197
$ws.Range("C133").Value = 45108.809027777803
$ws.Range("C134").Value = 45109.548611111102
$ws.Range("C135").Value = 45109.638888888898
$ws.Range("C136").Value = 45109.694444444402
$ws.Range("C137").Value = 45113.805555555598
$ws.Range("C138").Value = 45114.826388888898
$ws.Range("C139").Value = 45115.572916666701
$ws.Range("C140").Value = 45115.690972222197
$ws.Range("C141").Value = 45115.809027777803
$ws.Range("C142").Value = 45115.819444444402
$ws.Range("C143").Value = 45116.548611111102
$ws.Range("C144").Value = 45116.638888888898
$ws.Range("C145").Value = 45116.694444444402
$ws.Range("C146").Value = 45120.805555555598
$ws.Range("C147").Value = 45121.826388888898
$ws.Range("C148").Value = 45122.572916666701
$ws.Range("C149").Value = 45122.590277777803
$ws.Range("C150").Value = 45122.690972222197
$ws.Range("C151").Value = 45122.809027777803
$ws.Range("C152").Value = 45122.819444444402
$ws.Range("C153").Value = 45123.548611111102
$ws.Range("C154").Value = 45123.694444444402
$ws.Range("C155").Value = 45128.826388888898
$ws.Range("C156").Value = 45129.572916666701
$ws.Range("C157").Value = 45129.590277777803
$ws.Range("C158").Value = 45129.690972222197
$ws.Range("C159").Value = 45129.819444444402
$ws.Range("C160").Value = 45129.819444444402
$ws.Range("C161").Value = 45130.548611111102
$ws.Range("C162").Value = 45130.638888888898
$ws.Range("C163").Value = 45130.694444444402
$ws.Range("C164").Value = 45135.826388888898
$ws.Range("C165").Value = 45136.572916666701
$ws.Range("C166").Value = 45136.572916666701
$ws.Range("C167").Value = 45136.690972222197
$ws.Range("C168").Value = 45136.809027777803
$ws.Range("C169").Value = 45136.819444444402
$ws.Range("C170").Value = 45137.548611111102
$ws.Range("C171").Value = 45137.638888888898
$ws.Range("C172").Value = 45137.694444444402
$ws.Range("C173").Value = 45142.826388888898
$ws.Range("C174").Value = 45143.572916666701
$ws.Range("C175").Value = 45143.590277777803
$ws.Range("C176").Value = 45143.690972222197
$ws.Range("C177").Value = 45143.809027777803
$ws.Range("C178").Value = 45143.8125
$ws.Range("C179").Value = 45144.548611111102
$ws.Range("C180").Value = 45144.638888888898
$ws.Range("C181").Value = 45144.694444444402
$ws.Range("C182").Value = 45149.826388888898
$ws.Range("C183").Value = 45150.572916666701
$ws.Range("C184").Value = 45150.590277777803
$ws.Range("C185").Value = 45150.690972222197
$ws.Range("C186").Value = 45150.809027777803
$ws.Range("C187").Value = 45150.840277777803
$ws.Range("C188").Value = 45151.548611111102
$ws.Range("C189").Value = 45151.638888888898
$ws.Range("C190").Value = 45151.694444444402
$ws.Range("C191").Value = 45156.826388888898
$ws.Range("C192").Value = 45157.572916666701
$ws.Range("C193").Value = 45157.590277777803
$ws.Range("C194").Value = 45157.690972222197
$ws.Range("C195").Value = 45157.809027777803
$ws.Range("C196").Value = 45157.819444444402
$ws.Range("C197").Value = 45158.548611111102
$ws.Range("C198").Value = 45158.638888888898
$ws.Range("C199").Value = 45158.694444444402

# Restore the saved cursor/selection position from the source workbook.
$ws.Range("D196").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 175
    $win.ScrollColumn = 1
} catch {
    # Older/limited hosts may not expose window scroll position; selection
    # above already captures the meaningful, persisted state.
}
